# Apply updated market-board / profit values to the Ixion Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2394456.5
$ws.Range("J17").Value = 2394456.5
$ws.Range("L17").Value = 7183369.5
$ws.Range("N17").Value = -7183705.5
# Row 19
$ws.Range("H19").Value = 611.6667
$ws.Range("I19").Value = 600
$ws.Range("J19").Value = 626.25
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 626.25
$ws.Range("M19").Value = -425
$ws.Range("N19").Value = -976.25
# Row 21
$ws.Range("H21").Value = 33496.855
$ws.Range("I21").Value = 47759.5
$ws.Range("J21").Value = 14480
$ws.Range("K21").Value = 47759.5
$ws.Range("L21").Value = 14480
$ws.Range("M21").Value = -47291.5
$ws.Range("N21").Value = -15416
# Row 23
$ws.Range("H23").Value = 33496.855
$ws.Range("I23").Value = 47759.5
$ws.Range("J23").Value = 14480
$ws.Range("K23").Value = 47759.5
$ws.Range("L23").Value = 14480
$ws.Range("M23").Value = -47525.5
$ws.Range("N23").Value = -14948
# Row 53
$ws.Range("H53").Value = 59607.117
$ws.Range("I53").Value = 1571.2858
$ws.Range("J53").Value = 100232.2
$ws.Range("K53").Value = 1571.2858
$ws.Range("L53").Value = 100232.2
$ws.Range("M53").Value = -934.2858000000001
$ws.Range("N53").Value = -101506.2
# Row 112
$ws.Range("H112").Value = 5061.6416
$ws.Range("I112").Value = 663.3333
$ws.Range("J112").Value = 5325.54
$ws.Range("K112").Value = 1989.9999
$ws.Range("L112").Value = 15976.62
$ws.Range("M112").Value = -881.9999
$ws.Range("N112").Value = -18192.62
# Row 116
$ws.Range("H116").Value = 7145.75
$ws.Range("I116").Value = 10837.728
$ws.Range("K116").Value = 10837.728
$ws.Range("M116").Value = -7395.727999999999
# Row 137
$ws.Range("H137").Value = 1358.8077
$ws.Range("I137").Value = 1116.8572
$ws.Range("J137").Value = 2375
$ws.Range("K137").Value = 3350.5716
$ws.Range("L137").Value = 7125
$ws.Range("M137").Value = -800.5715999999998
$ws.Range("N137").Value = -12225
# Row 138
$ws.Range("H138").Value = 5645.478
$ws.Range("I138").Value = 980.7826
$ws.Range("J138").Value = 10310.174
$ws.Range("K138").Value = 2942.3478
$ws.Range("L138").Value = 30930.522
$ws.Range("M138").Value = 2197.6522
$ws.Range("N138").Value = -41210.522

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 4471.3438
$ws.Range("I45").Value = 6004.95
$ws.Range("J45").Value = 1915.3334
$ws.Range("K45").Value = 6004.95
$ws.Range("L45").Value = 1915.3334
$ws.Range("M45").Value = -5627.95
$ws.Range("N45").Value = -2669.3334
# Row 61
$ws.Range("H61").Value = 5909.72
$ws.Range("I61").Value = 6254.0435
$ws.Range("K61").Value = 6254.0435
$ws.Range("M61").Value = -6042.0435
# Row 68
$ws.Range("H68").Value = 61333.332
# Row 71
$ws.Range("H71").Value = 61333.332
# Row 74
$ws.Range("H74").Value = 1786.5454
$ws.Range("I74").Value = 1699.2157
$ws.Range("K74").Value = 1699.2157
$ws.Range("M74").Value = -825.2157
# Row 77
$ws.Range("H77").Value = 1786.5454
$ws.Range("I77").Value = 1699.2157
$ws.Range("K77").Value = 8496.0785
$ws.Range("M77").Value = -4128.0785
# Row 110
$ws.Range("H110").Value = 784.61536
$ws.Range("I110").Value = 759.70966
$ws.Range("J110").Value = 881.125
$ws.Range("K110").Value = 759.70966
$ws.Range("L110").Value = 881.125
$ws.Range("M110").Value = 1285.29034
$ws.Range("N110").Value = -4971.125
# Row 132
$ws.Range("H132").Value = 2770.7454
$ws.Range("I132").Value = 1527.08
$ws.Range("J132").Value = 3807.1333
$ws.Range("K132").Value = 4581.24
$ws.Range("L132").Value = 11421.3999
$ws.Range("M132").Value = -2051.24
$ws.Range("N132").Value = -16481.3999
# Row 136
$ws.Range("H136").Value = 5909.72
$ws.Range("I136").Value = 6254.0435
$ws.Range("K136").Value = 18762.1305
$ws.Range("M136").Value = -16212.1305

$ws = $wb.Worksheets.Item("BSM")
# Row 99
$ws.Range("H99").Value = 45455600
$ws.Range("I99").Value = 62501004
$ws.Range("J99").Value = 1185.1666
$ws.Range("K99").Value = 62501004
$ws.Range("L99").Value = 1185.1666
$ws.Range("M99").Value = -62499506
$ws.Range("N99").Value = -4181.1666
# Row 105
$ws.Range("H105").Value = 17665.691
$ws.Range("I105").Value = 29658.572
$ws.Range("J105").Value = 3674
$ws.Range("K105").Value = 29658.572
$ws.Range("L105").Value = 3674
$ws.Range("M105").Value = -27911.572
$ws.Range("N105").Value = -7168
# Row 107
$ws.Range("H107").Value = 753.0741
$ws.Range("I107").Value = 738.6957
$ws.Range("J107").Value = 835.75
$ws.Range("K107").Value = 738.6957
$ws.Range("L107").Value = 835.75
$ws.Range("M107").Value = 1181.3043
$ws.Range("N107").Value = -4675.75
# Row 134
$ws.Range("H134").Value = 3935.2856
$ws.Range("I134").Value = 5088.933
$ws.Range("J134").Value = 2113.7368
$ws.Range("K134").Value = 15266.799
$ws.Range("L134").Value = 6341.2104
$ws.Range("M134").Value = -12731.799
$ws.Range("N134").Value = -11411.2104

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4961.3945
$ws.Range("I31").Value = 1293.2222
$ws.Range("J31").Value = 13965.091
$ws.Range("K31").Value = 1293.2222
$ws.Range("L31").Value = 13965.091
$ws.Range("M31").Value = -998.2221999999999
$ws.Range("N31").Value = -14555.091
# Row 34
$ws.Range("H34").Value = 4961.3945
$ws.Range("I34").Value = 1293.2222
$ws.Range("J34").Value = 13965.091
$ws.Range("K34").Value = 1293.2222
$ws.Range("L34").Value = 13965.091
$ws.Range("M34").Value = -1091.2222
$ws.Range("N34").Value = -14369.091
# Row 58
$ws.Range("H58").Value = 1514.0444
$ws.Range("I58").Value = 739.913
$ws.Range("J58").Value = 2323.3635
$ws.Range("K58").Value = 739.913
$ws.Range("L58").Value = 2323.3635
$ws.Range("M58").Value = -536.913
$ws.Range("N58").Value = -2729.3635
# Row 81
$ws.Range("H81").Value = 33333.332
$ws.Range("J81").Value = 33333.332
$ws.Range("L81").Value = 33333.332
$ws.Range("N81").Value = -35329.332
# Row 84
$ws.Range("H84").Value = 33333.332
$ws.Range("J84").Value = 33333.332
$ws.Range("L84").Value = 99999.99600000001
$ws.Range("N84").Value = -109983.996
# Row 99
$ws.Range("H99").Value = 5002767
$ws.Range("I99").Value = 1674.3125
$ws.Range("J99").Value = 13893599
$ws.Range("K99").Value = 1674.3125
$ws.Range("L99").Value = 13893599
$ws.Range("M99").Value = -176.3125
$ws.Range("N99").Value = -13896595
# Row 122
$ws.Range("H122").Value = 1110.5454
$ws.Range("I122").Value = 913.82355
$ws.Range("K122").Value = 2741.47065
$ws.Range("M122").Value = -291.4706499999998
# Row 126
$ws.Range("H126").Value = 5002767
$ws.Range("I126").Value = 1674.3125
$ws.Range("J126").Value = 13893599
$ws.Range("K126").Value = 5022.9375
$ws.Range("L126").Value = 41680797
$ws.Range("M126").Value = -2552.9375
$ws.Range("N126").Value = -41685737
# Row 132
$ws.Range("H132").Value = 2369.7307
$ws.Range("I132").Value = 2017.6316
$ws.Range("J132").Value = 3325.4285
$ws.Range("K132").Value = 6052.8948
$ws.Range("L132").Value = 9976.2855
$ws.Range("M132").Value = -3522.8948
$ws.Range("N132").Value = -15036.2855
# Row 134
$ws.Range("H134").Value = 3291.4348
$ws.Range("I134").Value = 3793.4707
$ws.Range("J134").Value = 1869
$ws.Range("K134").Value = 11380.4121
$ws.Range("L134").Value = 5607
$ws.Range("M134").Value = -8845.4121
$ws.Range("N134").Value = -10677
# Row 136
$ws.Range("H136").Value = 1514.0444
$ws.Range("I136").Value = 739.913
$ws.Range("J136").Value = 2323.3635
$ws.Range("K136").Value = 2219.739
$ws.Range("L136").Value = 6970.0905
$ws.Range("M136").Value = 330.261
$ws.Range("N136").Value = -12070.0905

$ws = $wb.Worksheets.Item("CUL")
# Row 121
$ws.Range("H121").Value = 963.383
$ws.Range("I121").Value = 436.66666
$ws.Range("J121").Value = 1040.4634
$ws.Range("K121").Value = 1309.99998
$ws.Range("L121").Value = 3121.3902
$ws.Range("M121").Value = 0.00001999999994950485
$ws.Range("N121").Value = -5741.3902

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 5133.3667
$ws.Range("I126").Value = 5858.3335
$ws.Range("K126").Value = 17575.0005
$ws.Range("M126").Value = -15105.0005
# Row 132
$ws.Range("H132").Value = 3496.9062
$ws.Range("I132").Value = 4767.3076
$ws.Range("J132").Value = 2627.6843
$ws.Range("K132").Value = 14301.9228
$ws.Range("L132").Value = 7883.0529
$ws.Range("M132").Value = -11771.9228
$ws.Range("N132").Value = -12943.0529

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 44032.25
$ws.Range("I7").Value = 78959.16
$ws.Range("J7").Value = 2755
$ws.Range("K7").Value = 78959.16
$ws.Range("L7").Value = 2755
$ws.Range("M7").Value = -78847.16
$ws.Range("N7").Value = -2979
# Row 126
$ws.Range("H126").Value = 44032.25
$ws.Range("I126").Value = 78959.16
$ws.Range("J126").Value = 2755
$ws.Range("K126").Value = 236877.48
$ws.Range("L126").Value = 8265
$ws.Range("M126").Value = -234407.48
$ws.Range("N126").Value = -13205
# Row 132
$ws.Range("H132").Value = 25472274
$ws.Range("I132").Value = 30565668
$ws.Range("K132").Value = 91697004
$ws.Range("M132").Value = -91694474
# Row 136
$ws.Range("H136").Value = 4199.091
$ws.Range("I136").Value = 4569.5
$ws.Range("K136").Value = 13708.5
$ws.Range("M136").Value = -11158.5

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2077.76
$ws.Range("I122").Value = 1361.8334
$ws.Range("J122").Value = 2738.6155
$ws.Range("K122").Value = 4085.5002
$ws.Range("L122").Value = 8215.8465
$ws.Range("M122").Value = -1635.5002
$ws.Range("N122").Value = -13115.8465
# Row 126
$ws.Range("H126").Value = 1298
$ws.Range("I126").Value = 1298
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3894
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1424
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 2098.652
$ws.Range("I132").Value = 1475.6923
$ws.Range("J132").Value = 2908.5
$ws.Range("K132").Value = 4427.0769
$ws.Range("L132").Value = 8725.5
$ws.Range("M132").Value = -1897.0769
$ws.Range("N132").Value = -13785.5
# Row 136
$ws.Range("H136").Value = 4106.15
$ws.Range("I136").Value = 9082.5
$ws.Range("J136").Value = 1973.4286
$ws.Range("K136").Value = 27247.5
$ws.Range("L136").Value = 5920.2858
$ws.Range("M136").Value = -24697.5
$ws.Range("N136").Value = -11020.2858
